$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shared-string / label edits that happen to already-existing cells
# ---------------------------------------------------------------------------
# (handled further down once the row/col layout has been reshaped)

# ---------------------------------------------------------------------------
# 2. Reshape the grid: insert the extra "No" columns, drop the old legend
#    row, insert the new "Register Inputs" header row.
# ---------------------------------------------------------------------------
$ws.Columns("C:C").Insert()
$ws.Columns("F:F").Insert()
$ws.Rows("10:10").Delete()
$ws.Rows("2:2").Insert()

# ---------------------------------------------------------------------------
# 3. Column widths (B:K all 13.28515625)
# ---------------------------------------------------------------------------
$ws.Range("B1:K1").ColumnWidth = 13.28515625

# ---------------------------------------------------------------------------
# 4. Row 1 - function headers (merged groups)
# ---------------------------------------------------------------------------
$ws.Range("B1:D1").Merge()
$ws.Range("E1:G1").Merge()
$ws.Range("H1:I1").Merge()
$ws.Range("J1:K1").Merge()

$ws.Range("B1").Value = "Madgwick_qDot"
$ws.Range("E1").Value = "Madgwick_normalize"
$ws.Range("H1").Value = "Madgwick_correction"
$ws.Range("J1").Value = "Madgwick_update"

# ---------------------------------------------------------------------------
# 5. Row 2 - new "Register Inputs" No/Yes sub-header
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Register Inputs"
$ws.Range("B2").Value = "No"
$ws.Range("C2").Value = "Yes"
$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "Yes"
$ws.Range("G2").Value = "Yes"
$ws.Range("H2").Value = "No"
$ws.Range("I2").Value = "Yes"
$ws.Range("J2").Value = "No"
$ws.Range("K2").Value = "Yes"

# ---------------------------------------------------------------------------
# 6. Data rows - fill in the values for the newly inserted C & F columns,
#    and correct any other cells that moved/changed.
# ---------------------------------------------------------------------------
# Row 3 - Resource Sharing Factor
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4

# Row 4 - Multipliers
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 12
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 124
$ws.Range("I4").Value = 49
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 3

# Row 5 - Adders
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 133
$ws.Range("I5").Value = 133
$ws.Range("J5").Value = 8
$ws.Range("K5").Value = 8

# Row 6 - Registers
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 55
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 20
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 734
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 59

# Row 7 - RAMs
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0

# Row 8 - Multiplexers
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 35
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 28
$ws.Range("H8").Value = 12
$ws.Range("I8").Value = 361
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 24

# ---------------------------------------------------------------------------
# 7. Formatting
# ---------------------------------------------------------------------------
# Row 2 + the two "grouping" cells in row 3 (C3 & F3) get the "Note" look
# (tan/cream fill, thin grey box border) with centered text.
$noteRanges = @("A2:K2", "C3", "F3")
foreach ($addr in $noteRanges) {
    $rng = $ws.Range($addr)
    $rng.Style = "Note"
    $rng.Borders.Item(7).Color = 11711154
    $rng.Borders.Item(7).LineStyle = 1
    $rng.Borders.Item(8).Color = 11711154
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(9).Color = 11711154
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(10).Color = 11711154
    $rng.Borders.Item(10).LineStyle = 1
    $rng.HorizontalAlignment = -4108
}

# Row 1 merged header groups: centered "Note" style cells, with a grey
# underline (bottom border) spanning the whole group.
$headerGroups = @("B1:D1", "E1:G1", "H1:I1", "J1:K1")
foreach ($addr in $headerGroups) {
    $rng = $ws.Range($addr)
    $rng.Style = "Note"
    $rng.Borders.Item(9).Color = 11711154
    $rng.Borders.Item(9).LineStyle = 1
    $rng.HorizontalAlignment = -4108
}
# Left/right outline on the outer edges of each header group
$ws.Range("B1:D1").Borders.Item(7).Color = 11711154
$ws.Range("B1:D1").Borders.Item(7).LineStyle = 1
$ws.Range("B1:D1").Borders.Item(10).Color = 11711154
$ws.Range("B1:D1").Borders.Item(10).LineStyle = 1

$ws.Range("E1:G1").Borders.Item(7).Color = 11711154
$ws.Range("E1:G1").Borders.Item(7).LineStyle = 1

$ws.Range("H1:I1").Borders.Item(7).Color = 11711154
$ws.Range("H1:I1").Borders.Item(7).LineStyle = 1

$ws.Range("J1:K1").Borders.Item(7).Color = 11711154
$ws.Range("J1:K1").Borders.Item(7).LineStyle = 1

# ---------------------------------------------------------------------------
# 8. View: freeze panes below row 3 / right of col A, update selections
# ---------------------------------------------------------------------------
$ws.Range("B4").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("K6").Select()
